$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Text = "Możesz wprowadzić się w szał kosztem akcji dodatkowej. +4 do celności oraz -2 do Obrony i jesteś uciszony."
$find.Replacement.Text = "Atletyka."
$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)
